$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert B2's existing font so the engine re-resolves it onto the
# already-identical style used by B3 (Inter 12 FFF4594E) instead of
# keeping its own separate (duplicate) style entry.
$ws.Range("B2").Font.Name = "Inter"
$ws.Range("B2").Font.Size = 12
$ws.Range("B2").Font.Color = 5134836

# Update the test-data values for the "Custom Date" membership row.
$ws.Range("B2").Value = "Westwood Packages One Time"
$ws.Range("D2").Value = "8 September 2023"

# The old "Mosco Package 03" row is no longer part of the test data.
$ws.Rows("3:3").Delete()

# Column A keeps its existing width; column B grows to fit the new,
# longer package name.
$ws.Columns("B").ColumnWidth = 28.64

# Move the active selection, matching where the author left the cursor.
$ws.Range("F17").Select() | Out-Null
